# Updates cryptos list prices/volumes (Price column D, Volume(1h) column E)
# Commit message: Updated cryptos list on Fri Apr  7 17:46:20 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.010.53'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.862.68'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'" + '312.26'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = "'" + '0.5111'
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('D8').Value = "'" + '0.3863'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = "'" + '0.08279'
$ws.Range('E9').Value = '  -8.03%  '
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').Value = "'" + '41.60'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('E12').Value = '  -2.50%  '
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').Value = '1.864.87'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Value = "'" + '7.256'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = "'" + '90.66'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = "'" + '0.06653'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').Value = "'" + '17.71'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = "'" + '6.017'
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('D23').Value = '28.033.91'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = "'" + '11.09'
$ws.Range('E24').Value = '  -3.09%  '
$ws.Range('D25').Value = "'" + '2.234'
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').Value = '2.074.15'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').Value = "'" + '2.518'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = "'" + '157.64'
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('D30').Value = "'" + '124.89'
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('D31').Value = "'" + '0.1062'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('D33').Value = "'" + '5.962'
$ws.Range('E33').Value = '  +5.98%  '
$ws.Range('D34').Value = "'" + '3.588'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').Value = "'" + '9.398'
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').Value = "'" + '0.02413'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = "'" + '0.06493'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('D38').Value = "'" + '0.2176'
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('D39').Value = "'" + '0.6599'
$ws.Range('E39').Value = '  +2.80%  '
$ws.Range('D40').Value = "'" + '1.195'
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('D41').Value = "'" + '5.017'
$ws.Range('E41').Value = '  +1.75%  '
$ws.Range('D42').Value = "'" + '1.226'
$ws.Range('E42').Value = '  -4.70%  '
$ws.Range('E43').Value = '  -2.83%  '
$ws.Range('D44').Value = "'" + '0.6167'
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('D45').Value = "'" + '13.00'
$ws.Range('E45').Value = '  -1.70%  '
$ws.Range('D46').Value = "'" + '1.281'
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('D47').Value = "'" + '3.658'
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').Value = "'" + '1.213'
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').Value = "'" + '120.02'
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('D51').Value = "'" + '78.97'
$ws.Range('E51').Value = '  -0.83%  '
